# Fix "model" sheet: 'type' column should reflect the prompt/element type
# (e.g. geopoint instead of object, select_one instead of plain string) and
# the now-redundant 'elementType' column is dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")

# refrigerator_location: type was "object" (with elementType "geopoint") --
# now the type itself is "geopoint".
$ws.Cells.Item(4, 2).Value = "geopoint"

# refrigerator_condition: type was the generic "string" -- now it matches
# the prompt type used in the survey sheet, "select_one".
$ws.Cells.Item(5, 2).Value = "select_one"

# The elementType column (C) is no longer needed; remove it entirely.
$ws.Columns.Item(3).Delete()

# This sheet becomes the active / selected tab, with C22 as the last
# selection.
$ws.Activate()
$ws.Range("C22").Select()
